$wb = $excel.ActiveWorkbook

# Work on Jesse's tab: add the new time-log entry in row 6, mirroring the
# existing rows (date / minutes spent / aspect-of-project note).
$ws = $wb.Worksheets.Item("Jesse")

$ws.Cells.Item(6, 1).Value = "10/3102017"
$ws.Cells.Item(6, 2).Value = 120
$ws.Cells.Item(6, 3).Value = "Completed Room explicit constructor. Added input file for testing. "

# Row 6 wraps text like row 5 above it, so it gets the same taller row height.
$ws.Rows.Item(6).RowHeight = 28.5

# Leave the saved selection on Jesse at G11 (matching the author's last
# position there) without changing which sheet/tab is active overall.
$ws.Activate()
[void]$ws.Range("G11").Select()
$wb.Worksheets.Item("Main").Activate()

$wb.Save()
